$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 2: URL value
$ws.Range("B2").Value = "https://hl7.fr/fhir/fr/medication/ValueSet/FrMedicinalProductOnly"

# Row 8: Date value
$ws.Range("B8").Value = "2025-04-10T15:35:36+00:00"

# Row 14: Copyright value - clear it out entirely
$ws.Range("B14").ClearContents()
